# Update cryptocurrency price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.431.77'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '2.586.62'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.65'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.39'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '2.597.86'
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("E12").Value = '  +5.67%  '
$ws.Range("E13").Value = '  +4.85%  '
$ws.Range("D14").Value = '3.042.11'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '59.404.79'
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.96'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +5.32%  '
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '2.589.70'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.98'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("E21").Value = '  +2.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.57'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +6.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +8.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.97'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.48'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.08%  '
$ws.Range("D29").Value = '0.0₃0769'
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.10'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '157.64'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.34'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.07'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.917'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.06%  '
$ws.Range("E37").Value = '  +3.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.57'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.840'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.92%  '
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '289.47'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.60'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +9.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0974'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.597'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.65'
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0534'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("D50").Value = '1.975.38'
$ws.Range("E50").Value = '  +2.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.58'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.14%  '
